# Apply the diff to data.xlsx / Sheet1:
#  1. Fill in the previously-empty B:J data for rows 44-51 (daily metrics),
#     and the D-column "organic impressions" formula (B-C) for rows 44-61
#     (rows 52-61 have no B/C source data yet, so the formula evaluates to 0).
#  2. Update the frozen-pane view / active selection to D33 / H53.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1a. New row values for columns B..J, rows 44-51 -----------------------
$data = @{
    44 = @(3044, 2742, 69, 4, 6, 200, 8, 4649)
    45 = @(2671, 2310, 75, 2, 2, 157, 5, 70)
    46 = @(215,  4,    17, 0, 1, 17,  0, 898)
    47 = @(3572, 3260, 80, 2, 4, 160, 0, 0)
    48 = @(3380, 2995, 75, 9, 9, 284, 8, 5982)
    49 = @(4313, 3946, 83, 9, 8, 259, 9, 2658)
    50 = @(5737, 5271, 95, 8, 9, 294, 10, 89)
    51 = @(4179, 3765, 91, 3, 7, 235, 3, 941)
}

foreach ($row in 44..51) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("E$row").Value = $vals[2]
    $ws.Range("F$row").Value = $vals[3]
    $ws.Range("G$row").Value = $vals[4]
    $ws.Range("H$row").Value = $vals[5]
    $ws.Range("I$row").Value = $vals[6]
    $ws.Range("J$row").Value = $vals[7]
}

# --- 1b. D column "organic impressions" formula, rows 44-61 ----------------
# Assigning the formula across the whole contiguous range in one go lets the
# engine register it as a single shared formula group (same representation
# Excel itself uses for D3:D32 / D33:D43 already in this sheet).
$ws.Range("D44:D61").Formula = "=B44-C44"

# D32 used to carry its own standalone copy of the same formula as the
# D3:D32 shared group; bring it in line with the rest of that group.
$ws.Range("D4:D32").Formula = "=B4-C4"

# --- 2. Frozen-pane scroll position / active selection ---------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 33
$win.ScrollColumn = 4

$ws.Range("H53").Select()

Write-Output "edit applied"
